$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the daily log. It becomes the new row 7,
# pushing every existing record (old rows 7-41) down by one row (new rows 8-42).
$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value  = 11
$ws.Cells.Item(7, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value  = "Bíobío"
$ws.Cells.Item(7, 4).Value  = 44622
$ws.Cells.Item(7, 5).Value  = 8
$ws.Cells.Item(7, 6).Value  = 100112031
$ws.Cells.Item(7, 7).Value  = "Poroto verde"
$ws.Cells.Item(7, 8).Value  = "Magnum"
$ws.Cells.Item(7, 9).Value  = "Primera"
$ws.Cells.Item(7, 10).Value = 220
$ws.Cells.Item(7, 11).Value = 24000
$ws.Cells.Item(7, 12).Value = 25000
$ws.Cells.Item(7, 13).Value = 24545
$ws.Cells.Item(7, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(7, 15).Value = "Región del Maule"
$ws.Cells.Item(7, 16).Value = 982
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
